$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume (and any touched name/link) cells stay as text,
# not auto-converted to numbers by Excel when we assign their new values.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '71.621.73'
$ws.Range("E2").Value = '  +1.62%  '
$ws.Range("D3").Value = '3.817.00'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '703.09'
$ws.Range("E5").Value = '  +6.14%  '
$ws.Range("D6").Value = '175.23'
$ws.Range("E6").Value = '  +4.29%  '
$ws.Range("D7").Value = '3.815.32'
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '0.530'
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("D10").Value = '0.164'
$ws.Range("E10").Value = '  +1.98%  '
$ws.Range("D11").Value = '7.37'
$ws.Range("E11").Value = '  +5.87%  '
$ws.Range("D12").Value = '0.462'
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("E13").Value = '  +7.01%  '
$ws.Range("D14").Value = '36.61'
$ws.Range("E14").Value = '  +2.41%  '
$ws.Range("D15").Value = '4.457.72'
$ws.Range("E15").Value = '  +0.20%  '
$ws.Range("D16").Value = '3.798.12'
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("D17").Value = '71.548.36'
$ws.Range("E17").Value = '  +1.59%  '
$ws.Range("D18").Value = '17.75'
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").Value = '7.24'
$ws.Range("E19").Value = '  +1.19%  '
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").Value = '11.07'
$ws.Range("E21").Value = '  +6.05%  '
$ws.Range("D22").Value = '484.44'
$ws.Range("E22").Value = '  +2.04%  '
$ws.Range("D23").Value = '0.719'
$ws.Range("E23").Value = '  +0.69%  '
$ws.Range("D24").Value = '84.58'
$ws.Range("E24").Value = '  +2.25%  '
$ws.Range("E25").Value = '  -2.21%  '
$ws.Range("D26").Value = '12.36'
$ws.Range("E26").Value = '  +1.31%  '
$ws.Range("E27").Value = '  +1.92%  '
$ws.Range("E28").Value = '  +1.64%  '
$ws.Range("D29").Value = '3.966.92'
$ws.Range("E29").Value = '  +0.22%  '
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("E31").Value = '  +11.34%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '2.32'
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").Value = '7.63'
$ws.Range("E33").Value = '  +3.15%  '
$ws.Range("D34").Value = '0.189'
$ws.Range("E34").Value = '  +5.99%  '
$ws.Range("D35").Value = '29.68'
$ws.Range("E35").Value = '  +1.08%  '
$ws.Range("D36").Value = '9.32'
$ws.Range("E36").Value = '  +2.94%  '
$ws.Range("D38").Value = '0.104'
$ws.Range("E38").Value = '  +1.82%  '
$ws.Range("D39").Value = '3.46'
$ws.Range("E39").Value = '  +2.98%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '2.32'
$ws.Range("E40").Value = '  +11.00%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '6.04'
$ws.Range("E41").Value = '  +1.97%  '
$ws.Range("D42").Value = '0.993'
$ws.Range("E42").Value = '  +2.07%  '
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").Value = '164.58'
$ws.Range("E45").Value = '  +4.13%  '
$ws.Range("B46").Value = 'FLOKI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D46").Value = '0.000312'
$ws.Range("E46").Value = '  +11.67%  '
$ws.Range("D47").Value = '44.73'
$ws.Range("E47").Value = '  -0.93%  '
$ws.Range("D48").Value = '48.76'
$ws.Range("E48").Value = '  +2.02%  '
$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D49").Value = '420.53'
$ws.Range("E49").Value = '  +7.40%  '
$ws.Range("D50").Value = '0.304'
$ws.Range("E50").Value = '  +1.15%  '
$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").Value = '1.40'
$ws.Range("E51").Value = '  -1.82%  '

# Restore the default cell style so no stray number-format styling remains.
$ws.Range("B2:E51").Style = "Normal"
